{"js": "// Lab 6 proposal edit: broaden the hardware list and append a new sentence\n// describing the color-switching feature.\n//\n// Original sentence:\n//   \"For our Lab 6 project we are going to do the paint project with the\n//    camera, a keyboard, and a mouse.\"\n//\n// Target sentence:\n//   \"For our Lab 6 project we are going to do the paint project with the\n//    camera, a mouse, the VGA, and the speaker with a microphone.  The user\n//    will be able to change the colors that they can draw on the image with\n//    different colors, and switch between colors to draw with.\"\n\nconst body = context.document.body;\n\n// 1. Replace \"camera, a keyboard, and a mouse\" with the new, longer hardware\n//    list (\"camera, a mouse, the VGA, and the speaker with a microphone\").\nconst oldClause = body.search(\"camera, a keyboard, and a mouse\", { matchCase: true });\noldClause.load(\"items\");\nawait context.sync();\n\nif (oldClause.items.length > 0) {\n  oldClause.items[0].insertText(\n    \"camera, a mouse, the VGA, and the speaker with a microphone\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2. Insert the new sentence right after the sentence-ending period (and\n//    before the trailing _GoBack bookmark), so it lands inside the same\n//    paragraph.\nconst endClause = body.search(\"speaker with a microphone.\", { matchCase: true });\nendClause.load(\"items\");\nawait context.sync();\n\nif (endClause.items.length > 0) {\n  endClause.items[0].insertText(\n    \"  The user will be able to change the colors that they can draw on the image with different colors, and switch between colors to draw with.\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n", "ps1": "# Lab 6 proposal edit: broaden the hardware list and append a new sentence\n# describing the color-switching feature.\n#\n# Original sentence:\n#   \"For our Lab 6 project we are going to do the paint project with the\n#    camera, a keyboard, and a mouse.\"\n#\n# Target sentence:\n#   \"For our Lab 6 project we are going to do the paint project with the\n#    camera, a mouse, the VGA, and the speaker with a microphone.  The user\n#    will be able to change the colors that they can draw on the image with\n#    different colors, and switch between colors to draw with.\"\n\n$d = $word.ActiveDocument\n\n# 1. Replace \"camera, a keyboard, and a mouse\" with the new, longer hardware\n#    list (\"camera, a mouse, the VGA, and the speaker with a microphone\").\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"camera, a keyboard, and a mouse\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"camera, a mouse, the VGA, and the speaker with a microphone\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# 2. Insert the new sentence right after the sentence-ending period (and\n#    before the trailing _GoBack bookmark), so it lands inside the same\n#    paragraph.\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.ClearFormatting()\n$find2.Text = \"speaker with a microphone.\"\n$find2.Execute() | Out-Null\n$range2.Collapse(0)  # wdCollapseEnd\n$range2.InsertAfter(\"  The user will be able to change the colors that they can draw on the image with different colors, and switch between colors to draw with.\")\n"}
